$d = $word.ActiveDocument

$find = "The MANAGER has full control of other user profiles including the ability to change roles."
$replace = "The MANAGER has full control of other user profiles including the ability to change roles and assign employee and/or keyholder status."

$d.Content.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
